# Add working set of sequences
# For a fixed set of rows that currently only have data through column F,
# fill columns G through N with the same "N/A" value already used in F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4,6,9,13,17,20,27,30,32,34,35,38,41,43,44,46,49,52,53,55,59,60,61,65,68,71,74,78,82,85,89,93,97,100,102,107,108,109,116,121,126,128,131,132,133,136,140,144,150,152,156,158,162,164,168,169,173,178,184,187,188,190,191,193)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
